$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that follows the title heading.
$d.Paragraphs.Item(2).Range.Delete()

# 2. Insert a new bold paragraph ("Play Age of the Gods Norse: Ways of Thunder
#    for Free | Review") right before the final paragraph.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore()

$newParaRange = $d.Paragraphs.Item($count).Range
$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Age of the Gods Norse: Ways of Thunder for Free | Review</w:t></w:r></w:p>'
$newParaRange.InsertXML($newXml)

# 3. Replace the text of the (now last) italic paragraph with the new meta
#    description copy, keeping its existing italic formatting.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalPara.Range.Find.Execute(
    "Create a cartoon-style feature image for " + [char]34 + "Age of the Gods Norse Ways of Thunder" + [char]34 + ". The image should feature a happy Maya warrior with glasses. The Maya warrior should be holding a lightning bolt with Norse symbols on it, and standing in front of the reel set with diamond grid and blue electric coverage. The background should reflect the Norse mythology with a glimpse of the Norse gods and a thunderstorm. The image should attract the attention of online slot game enthusiasts and evoke the excitement of playing an epic adventure.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Age of the Gods Norse: Ways of Thunder and play for free. Dynamic gameplay, epic graphics, and exciting bonus features.",
    2)
